# Updates cryptos list values (Price and Volume(1h) columns), plus some
# coin name/link re-shuffles, per the commit "Updated cryptos list on
# Fri Aug 11 09:21:10 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.372.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'1.845.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'240.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'0.6300"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.2956"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "'24.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "'0.07721"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'1.853.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("D13").Value = "'4.995"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "'0.6835"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'0.000009990"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "'82.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "'6.138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "'29.411.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'227.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").Value = "'12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'7.547"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'157.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "'8.363"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'17.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.466"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.257"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05690"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("D31").Value = "'4.126"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "'4.013"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").Value = "'0.7146"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").Value = "'2.592"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'1.255.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'0.9129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "'6.206"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'2.017.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.31%  "
$ws.Range("D44").Value = "'100.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4023"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.118"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.693"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1126"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.10%  "
